$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new mapping rows (column A = FullyQualifiedName, column B = SqlName)
$ws.Range("A2").Value = "northwind.entities.City"
$ws.Range("B2").Value = "CITY"
$ws.Range("A3").Value = "northwind.entities.Territory#shipper"
$ws.Range("B3").Value = "ID_SHIPPER"
$ws.Range("A4").Value = "northwind.entities.Person#firstName"
$ws.Range("B4").Value = "FIRST_NAME"

# Re-apply the (unchanged-looking) font to the whole used range so a dedicated
# style/font entry is materialized and stamped on every cell, matching the
# original author's save.
$ws.Range("A1:B4").Font.ThemeColor = 1

# Page setup, as present in the saved workbook.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the cursor where the author's last interaction left it.
$result = $ws.Range("B14").Select()
